# Applies the cryptos-list price/volume refresh described in the commit
# 'Updated cryptos list on Fri Mar 22 05:35:00 UTC 2024 with GitHub Actions'.
#
# Numeric-looking values (e.g. "573.02") are written with a leading single
# quote so Excel keeps them as literal text instead of silently parsing them
# into numbers (which would both change the cell type and round/alter the
# displayed digits, e.g. "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $address, $text) {
    $sheet.Range($address).Value = $text
}

# Row 2
Set-TextCell $ws "D2" '66.180.59'
Set-TextCell $ws "E2" '  -0.59%  '

# Row 3
Set-TextCell $ws "D3" '3.512.77'
Set-TextCell $ws "E3" '  +1.16%  '

# Row 4
Set-TextCell $ws "D4" '''0.999'
Set-TextCell $ws "E4" '  +0.07%  '

# Row 5
Set-TextCell $ws "D5" '''573.02'
Set-TextCell $ws "E5" '  +5.26%  '

# Row 6
Set-TextCell $ws "D6" '''178.80'
Set-TextCell $ws "E6" '  -4.41%  '

# Row 7
Set-TextCell $ws "D7" '''0.638'
Set-TextCell $ws "E7" '  +5.82%  '

# Row 8
Set-TextCell $ws "D8" '''1.00'
Set-TextCell $ws "E8" '  +0.10%  '

# Row 9
Set-TextCell $ws "D9" '''0.635'
Set-TextCell $ws "E9" '  +1.47%  '

# Row 10
Set-TextCell $ws "E10" '  +4.75%  '

# Row 11
Set-TextCell $ws "D11" '''55.64'
Set-TextCell $ws "E11" '  +2.45%  '

# Row 12
Set-TextCell $ws "D12" '''0.0000275'
Set-TextCell $ws "E12" '  +3.38%  '

# Row 13
Set-TextCell $ws "D13" '''9.30'
Set-TextCell $ws "E13" '  +0.21%  '

# Row 14
Set-TextCell $ws "D14" '4.074.84'
Set-TextCell $ws "E14" '  +1.37%  '

# Row 15
Set-TextCell $ws "D15" '3.510.89'
Set-TextCell $ws "E15" '  +1.50%  '

# Row 16
Set-TextCell $ws "E16" '  +0.13%  '

# Row 17
Set-TextCell $ws "D17" '''18.38'
Set-TextCell $ws "E17" '  +2.14%  '

# Row 18
Set-TextCell $ws "D18" '66.175.93'

# Row 19
Set-TextCell $ws "D19" '''12.01'
Set-TextCell $ws "E19" '  +2.85%  '

# Row 20
Set-TextCell $ws "E20" '  +2.19%  '

# Row 21
Set-TextCell $ws "D21" '''414.95'
Set-TextCell $ws "E21" '  +0.32%  '

# Row 22
Set-TextCell $ws "D22" '''4.19'
Set-TextCell $ws "E22" '  +8.48%  '

# Row 23
Set-TextCell $ws "D23" '''4.29'
Set-TextCell $ws "E23" '  +3.40%  '

# Row 24
Set-TextCell $ws "D24" '''85.74'
Set-TextCell $ws "E24" '  +1.90%  '

# Row 25
Set-TextCell $ws "D25" '''13.11'
Set-TextCell $ws "E25" '  +11.56%  '

# Row 26
Set-TextCell $ws "D26" '''11.01'
Set-TextCell $ws "E26" '  -0.81%  '

# Row 27
Set-TextCell $ws "E27" '  -0.60%  '

# Row 28
Set-TextCell $ws "E28" '  +4.64%  '

# Row 29
Set-TextCell $ws "D29" '''30.49'
Set-TextCell $ws "E29" '  +2.20%  '

# Row 30
Set-TextCell $ws "D30" '''633.98'
Set-TextCell $ws "E30" '  -2.58%  '

# Row 31
Set-TextCell $ws "D31" '''6.51'
Set-TextCell $ws "E31" '  -1.25%  '

# Row 32
Set-TextCell $ws "D32" '''11.70'
Set-TextCell $ws "E32" '  +0.84%  '

# Row 33
Set-TextCell $ws "D33" '''0.111'
Set-TextCell $ws "E33" '  +1.28%  '

# Row 34
Set-TextCell $ws "E34" '  +14.19%  '

# Row 35
Set-TextCell $ws "D35" '''59.52'
Set-TextCell $ws "E35" '  +0.59%  '

# Row 36
Set-TextCell $ws "E36" '  +0.15%  '

# Row 37
Set-TextCell $ws "D37" '0.0₃0800'
Set-TextCell $ws "E37" '  -0.24%  '

# Row 38
Set-TextCell $ws "D38" '''37.27'
Set-TextCell $ws "E38" '  -2.31%  '

# Row 39
Set-TextCell $ws "D39" '''0.381'
Set-TextCell $ws "E39" '  -1.29%  '

# Row 40
Set-TextCell $ws "D40" '3.255.69'
Set-TextCell $ws "E40" '  +9.13%  '

# Row 41
Set-TextCell $ws "E41" '  +2.17%  '

# Row 42
Set-TextCell $ws "D42" '''0.999'
Set-TextCell $ws "E42" '  +0.12%  '

# Row 43
Set-TextCell $ws "E43" '  +1.78%  '

# Row 44
Set-TextCell $ws "D44" '''0.0420'
Set-TextCell $ws "E44" '  +2.06%  '

# Row 45
Set-TextCell $ws "B45" 'Fetch.AI'
Set-TextCell $ws "C45" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws "D45" '''2.53'
Set-TextCell $ws "E45" '  -3.45%  '

# Row 46
Set-TextCell $ws "B46" 'ApeXProtocol'
Set-TextCell $ws "C46" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws "D46" '''3.26'
Set-TextCell $ws "E46" '  -4.56%  '

# Row 47
Set-TextCell $ws "E47" '  +1.39%  '

# Row 48
Set-TextCell $ws "D48" '''0.133'
Set-TextCell $ws "E48" '  +2.91%  '

# Row 49
Set-TextCell $ws "D49" '''8.68'
Set-TextCell $ws "E49" '  -0.58%  '

# Row 50
Set-TextCell $ws "D50" '''138.59'
Set-TextCell $ws "E50" '  -0.37%  '

# Row 51
Set-TextCell $ws "D51" '''2.39'
Set-TextCell $ws "E51" '  -0.07%  '
